$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").Value = "Completed"
$ws.Range("E20").Value = "Incomplete"

$ws.Range("E4").Select()
